$wb = $excel.ActiveWorkbook

# --- Create "week 50" as a copy of "week 49" (keeps formatting/col widths/shared formula) ---
$week49 = $wb.Worksheets.Item("week 49")
$week49.Copy($null, $week49)
$week50 = $wb.Worksheets.Item("week 49 (2)")
$week50.Name = "week 50"

# --- Overwrite the row-7 entry with week 50's own log entry ---
$week50.Range("B7").Value = 41620
$week50.Range("C7").Value = 0.37152777777777773
$week50.Range("D7").Value = 0.38541666666666669
$week50.Range("F7").Value = "Alles bijgewerkt naar commit leraar"

# --- Clear the remaining (unused) detail rows copied from week 49 ---
$week50.Range("C8:D8").ClearContents()
$week50.Range("F8").ClearContents()
$week50.Range("C9:D9").ClearContents()
$week50.Range("F9").ClearContents()
$week50.Range("C10:D10").ClearContents()
$week50.Range("F10").ClearContents()
$week50.Range("C11:D11").ClearContents()
$week50.Range("F11").ClearContents()

# --- Selection / active sheet lands on week 50, cell F8 ---
$week50.Range("F8").Select()

# --- Update "Totaal" sheet: insert a row for week 50 before the summary row ---
$totaal = $wb.Worksheets.Item("Totaal")
$totaal.Rows.Item(9).Insert()

$totaal.Range("A8").Value = 49
$totaal.Range("B8").Formula = "='week 49'!G18"

$totaal.Range("A9").Value = 50
$totaal.Range("B9").Formula = "='week 50'!G18"

$totaal.Range("B10").Select()

# --- Re-select week 50 / F8 so it ends up the active tab (matches activeTab=2) ---
$week50.Range("F8").Select()
